$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values after repulling data / recalculating mean
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("F32").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("F45").Value = 2
$ws.Range("F53").Value = 4
$ws.Range("F66").Value = -5
